# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the a5b06ae6-1a02-44f6-88af-661e12ff4cda.md row (row 6) on each sheet to
# reflect a fresh handoff report being generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G is "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-21 16:49:44"

# zh-cn sheet: column H is "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-21 16:49:40"

# de-de sheet: column H is "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-21 16:49:44"
